# The deck currently uses the "Integral" design theme (green palette) for
# its slide master / slides. This edit switches the presentation's active
# colour scheme back to the standard default "Office" palette, which is
# what a user does by picking the first ("Office") swatch in
# Design > Variants > Colors (PowerPoint then rewrites the theme part that
# the slide master points at with the new scheme).
#
# ThemeColorScheme.Colors(i).RGB uses the VBA/OLE colour order (0xBBGGRR),
# so each target "RRGGBB" hex value below is byte-reversed before being
# assigned.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# index -> target theme colour (standard Office palette), in RRGGBB
$officeColors = @{
    1  = 0x000000   # dk1 / Background 1
    2  = 0xFFFFFF   # lt1 / Text 1
    3  = 0x44546A   # dk2 / Text 2
    4  = 0xE7E6E6   # lt2 / Background 2
    5  = 0x5B9BD5   # accent1
    6  = 0xED7D31   # accent2
    7  = 0xA5A5A5   # accent3
    8  = 0xFFC000   # accent4
    9  = 0x4472C4   # accent5
    10 = 0x70AD47   # accent6
    11 = 0x0563C1   # hyperlink
    12 = 0x954F72   # followed hyperlink
}

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $rgbHex = $officeColors[$i]
    $r = [math]::Floor($rgbHex / 0x10000) % 0x100
    $g = [math]::Floor($rgbHex / 0x100) % 0x100
    $b = $rgbHex % 0x100
    $bgr = ($b * 0x10000) + ($g * 0x100) + $r
    $colorScheme.Colors($i).RGB = $bgr
}
